# "quick fix in Class Diagram"
#
# The commit refreshes the stale date / slide-number footer fields on every
# slide (re-cached after reopening the deck later) and gives the
# "Preferences" class box a dashed outline on the Class Diagram slide.

$p = $ppt.ActivePresentation

# Refresh any date / slide-number field placeholders, on every slide.
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "23/05/15") {
                $tr.Text = "8-7-2015"
            } elseif ($tr.Text -eq [char]0x2039 + "nr." + [char]0x203A) {
                $tr.Text = [char]0x2039 + "#" + [char]0x203A
            }
        }
    }
}

# Give the "Preferences" class box (roundRect with a plain, invisible line)
# a visible dashed outline: <a:ln w="25400"><a:prstDash val="dash"/></a:ln>
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)
        if ($shp.Name -eq "Afgeronde rechthoek 21") {
            $shp.Line.Weight = 2
            $shp.Line.DashStyle = 4
        }
    }
}
